# Auto-generated Excel COM-interop script
# Commit: Updated cryptos list on Mon Nov 13 06:31:32 UTC 2023 with GitHub Actions
# Applies the Price (D) / Volume(1h) (E) text updates for rows 2-51 of the cryptos sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.890.22'
$ws.Range('E2').Value = '  -0.64%  '
$ws.Range('D3').Value = '2.034.95'
$ws.Range('E3').Value = '  -0.90%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.92'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.82%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.654'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.45%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '58.02'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.13%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.374'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0767'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.11%  '
$ws.Range('E11').Value = '  +2.12%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.24'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -5.04%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.875'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +7.57%  '
$ws.Range('D14').Value = '2.333.50'
$ws.Range('E14').Value = '  -0.81%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.58'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.97%  '
$ws.Range('D16').Value = '2.039.35'
$ws.Range('E16').Value = '  -0.75%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '18.02'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +6.09%  '
$ws.Range('D18').Value = '36.878.89'
$ws.Range('E18').Value = '  -0.71%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '73.34'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.68%  '
$ws.Range('D20').Value = '0.0₃0881'
$ws.Range('E20').Value = '  -2.29%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.32'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.73%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '234.50'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.87%  '
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.45'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.19%  '
$ws.Range('E25').Value = '  +3.00%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '167.57'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.41%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.13'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.70%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.80'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.50%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.42'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +14.78%  '
$ws.Range('E30').Value = '  -0.82%  '
$ws.Range('E31').Value = '  -4.46%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.66'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.55%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0609'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.87%  '
$ws.Range('E34').Value = '  +0.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0865'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.30%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.83'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.65%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.22'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.42%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.29'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.70%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.11'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.28%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.17'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.21%  '
$ws.Range('E41').Value = '  -0.16%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0948'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -12.78%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.13'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.39%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '96.43'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.39%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '16.76'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.38%  '
$ws.Range('D46').Value = '1.286.94'
$ws.Range('E46').Value = '  +0.35%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.34'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.42%  '
$ws.Range('E48').Value = '  -0.37%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.73'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +6.26%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.66'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.84%  '
$ws.Range('D51').Value = '2.220.74'
$ws.Range('E51').Value = '  -0.94%  '
